$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.758.49'
$ws.Range('E2').Value = '  -4.32%  '
$ws.Range('D3').Value = '2.452.51'
$ws.Range('E3').Value = '  -5.66%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''548.28'
$ws.Range('E5').Value = '  -3.95%  '
$ws.Range('D6').Value = '''144.90'
$ws.Range('E6').Value = '  -5.97%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.595'
$ws.Range('E8').Value = '  -3.64%  '
$ws.Range('D9').Value = '2.448.64'
$ws.Range('E9').Value = '  -5.70%  '
$ws.Range('E10').Value = '  -7.50%  '
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').Value = '''5.38'
$ws.Range('E12').Value = '  -6.99%  '
$ws.Range('E13').Value = '  -6.84%  '
$ws.Range('D14').Value = '''25.96'
$ws.Range('E14').Value = '  -6.55%  '
$ws.Range('D15').Value = '2.893.44'
$ws.Range('E15').Value = '  -5.71%  '
$ws.Range('E16').Value = '  -8.09%  '
$ws.Range('D17').Value = '60.675.37'
$ws.Range('E17').Value = '  -4.26%  '
$ws.Range('D18').Value = '2.453.78'
$ws.Range('E18').Value = '  -5.50%  '
$ws.Range('D19').Value = '''11.03'
$ws.Range('E19').Value = '  -7.14%  '
$ws.Range('D20').Value = '''6.89'
$ws.Range('E20').Value = '  -7.36%  '
$ws.Range('E21').Value = '  -6.75%  '
$ws.Range('D22').Value = '''318.44'
$ws.Range('E22').Value = '  -6.02%  '
$ws.Range('D24').Value = '''63.31'
$ws.Range('E24').Value = '  -5.60%  '
$ws.Range('D25').Value = '''1.76'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').Value = '0.0₃0978'
$ws.Range('E26').Value = '  -6.42%  '
$ws.Range('D27').Value = '2.574.75'
$ws.Range('E27').Value = '  -5.76%  '
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').Value = '''1.49'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = '''535.24'
$ws.Range('E30').Value = '  -7.26%  '
$ws.Range('D31').Value = '''8.32'
$ws.Range('E31').Value = '  -7.93%  '
$ws.Range('D32').Value = '''7.58'
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('E33').Value = '  -6.42%  '
$ws.Range('D34').Value = '''1.89'
$ws.Range('E34').Value = '  -7.16%  '
$ws.Range('D35').Value = '''1.57'
$ws.Range('E35').Value = '  -7.76%  '
$ws.Range('D36').Value = '''5.84'
$ws.Range('E36').Value = '  -10.33%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').Value = '''4.82'
$ws.Range('E38').Value = '  -9.15%  '
$ws.Range('D39').Value = '''0.376'
$ws.Range('E39').Value = '  -5.64%  '
$ws.Range('D40').Value = '''18.40'
$ws.Range('E40').Value = '  -5.69%  '
$ws.Range('D41').Value = '''144.30'
$ws.Range('E41').Value = '  -6.24%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('E43').Value = '  -8.24%  '
$ws.Range('D44').Value = '''39.71'
$ws.Range('E44').Value = '  -4.30%  '
$ws.Range('D45').Value = '''2.28'
$ws.Range('E45').Value = '  -7.97%  '
$ws.Range('D46').Value = '''146.33'
$ws.Range('E46').Value = '  -6.49%  '
$ws.Range('E47').Value = '  -7.13%  '
$ws.Range('D48').Value = '''20.79'
$ws.Range('E48').Value = '  -9.91%  '
$ws.Range('D49').Value = '''0.0528'
$ws.Range('E49').Value = '  -8.58%  '
$ws.Range('D50').Value = '''0.581'
$ws.Range('E50').Value = '  -7.01%  '
$ws.Range('D51').Value = '''0.0936'
$ws.Range('E51').Value = '  -5.67%  '
